# Append a new data row (row 7) to the survey sheet, mirroring a fresh
# "upload data on google sheet" export: every cell in the new row is
# stored as literal text (numeric-looking values included), matching the
# existing rows' convention (e.g. row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A7:I7")

# Pre-format the new row as Text so that numeric-looking strings ("10",
# "5000", dates, …) are stored verbatim instead of being coerced into
# numbers - consistent with how this sheet already stores row 6.
$rng.NumberFormat = "@"

$ws.Range("A7").Value = "nettoyant"
$ws.Range("B7").Value = "10"
$ws.Range("C7").Value = "2025-04-04"
$ws.Range("D7").Value = "5000"
$ws.Range("E7").Value = "6000"
$ws.Range("F7").Value = "réduction"
$ws.Range("G7").Value = "homme"
$ws.Range("H7").Value = "25"
$ws.Range("I7").Value = "5000"

# Drop back to the default "Normal" style so we don't leave a stray
# number-format override applied to the new cells.
$rng.Style = "Normal"

# Extend the "numbers stored as text" ignored-error suppression to cover
# the newly-added row, same as the rest of the imported table.
try {
    $ws.Range("A1:I7").Errors.Item(1).Ignore = $true
} catch {
}
